$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New summary block appended below the existing tables (rows 51-64).
# Mirrors the "Asians / Black / White" arrest-count breakdown for 2018 vs
# 2019 plus the year-on-year percentage-change figures.
# ---------------------------------------------------------------------------

# Row 51 - headers
$ws.Range("D51").Value = 2018
$ws.Range("E51").Value = "Asians"
$ws.Range("F51").Value = "Black"
$ws.Range("G51").Value = "White"
$ws.Range("I51").Value = 2019
$ws.Range("J51").Value = "Asians"
$ws.Range("K51").Value = "Black"
$ws.Range("L51").Value = "White"

# Row 52
$ws.Range("E52").Value = 472
$ws.Range("F52").Value = 1250
$ws.Range("G52").Value = 1255
$ws.Range("J52").Value = 571
$ws.Range("K52").Value = 1351
$ws.Range("L52").Value = 1483
$ws.Range("L52").NumberFormat = "#,##0"
$ws.Range("L52").HorizontalAlignment = -4108
$ws.Range("L52").VerticalAlignment = -4108

# Row 53
$ws.Range("E53").Value = 569
$ws.Range("F53").Value = 1191
$ws.Range("G53").Value = 2815
$ws.Range("J53").Value = 576
$ws.Range("K53").Value = 1319
$ws.Range("L53").Value = 3041
$ws.Range("L53").NumberFormat = "#,##0"
$ws.Range("L53").HorizontalAlignment = -4108
$ws.Range("L53").VerticalAlignment = -4108

# Row 54
$ws.Range("E54").Value = 254
$ws.Range("F54").Value = 1167
$ws.Range("G54").Value = 147
$ws.Range("J54").Value = 289
$ws.Range("K54").Value = 1122
$ws.Range("L54").Value = 133
$ws.Range("L54").NumberFormat = "#,##0"
$ws.Range("L54").HorizontalAlignment = -4108
$ws.Range("L54").VerticalAlignment = -4108

# Row 55
$ws.Range("E55").Value = 235
$ws.Range("J55").Value = 246

# Row 56 - totals
$ws.Range("E56").Formula = "=SUM(E52:E55)"
$ws.Range("F56").Formula = "=SUM(F52:F55)"
$ws.Range("G56").Formula = "=SUM(G52:G55)"
$ws.Range("J56").Formula = "=SUM(J52:J55)"
$ws.Range("K56").Formula = "=SUM(K52:K55)"
$ws.Range("L56").Formula = "=SUM(L52:L55)"
$ws.Range("L56").NumberFormat = "#,##0"

# Rows 58-61 - helper column O (2019 Asian/Black/White/no-age-provided figures)
$ws.Range("O58").Value = 571
$ws.Range("O59").Value = 576
$ws.Range("O60").Value = 289
$ws.Range("O61").Value = 246

# Row 62 - White % change
$ws.Range("E62").Value = "W"
$ws.Range("F62").Formula = "=(4657-4217)/4217*100"
$ws.Range("O62").Formula = "=SUM(O58:O61)"

# Row 63 - Asian % change
$ws.Range("E63").Value = "A"
$ws.Range("F63").Formula = "=(1682-1530)/1530*100"

# Row 64 - Black % change
$ws.Range("E64").Value = "B"
$ws.Range("F64").Formula = "=(3792-3608)/3608*100"

# ---------------------------------------------------------------------------
# Scroll the view down to the new block and select the next empty cell,
# matching where the author's cursor ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("F65").Select()
